$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-33 from 45204 (2023-10-05)
# to 45207 (2023-10-08), matching the author's bulk date-shift update.
$ws.Range("C2:C33").Value = 45207
